$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update the date serial value in A1 (merged A1:D1) ---
$ws.Range("A1").Value = 45311

# --- Update prices in column D (rows 34-36) ---
$ws.Range("D34").Value = 100.382
$ws.Range("D35").Value = 142.229
$ws.Range("D36").Value = 192.417
